# Applies the "Updated cryptos list" data refresh: new Price (column D)
# and Volume(1h) (column E) values for the coin rows in the active sheet.
#
# Column D sometimes holds plain decimal-looking text (e.g. "0.999"); it is
# written through a temporary Text number format so Excel keeps it as a
# string instead of parsing it into a floating point number, then the cell
# style is restored to Normal so no formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.783.13'
$ws.Range("E2").Value = '  -0.19%  '

# Row 3
$ws.Range("D3").Value = '3.823.12'
$ws.Range("E3").Value = '  +2.11%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '

# Row 7
$ws.Range("D7").Value = '3.822.04'
$ws.Range("E7").Value = '  +2.15%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.48%  '

# Row 10
$ws.Range("E10").Value = '  -0.59%  '

# Row 11
$ws.Range("E11").Value = '  +2.82%  '

# Row 12
$ws.Range("E12").Value = '  -0.32%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.28%  '

# Row 14
$ws.Range("E14").Value = '  -1.32%  '

# Row 15
$ws.Range("D15").Value = '4.449.58'
$ws.Range("E15").Value = '  +2.08%  '

# Row 16
$ws.Range("D16").Value = '3.816.19'
$ws.Range("E16").Value = '  +1.58%  '

# Row 17
$ws.Range("D17").Value = '69.816.69'
$ws.Range("E17").Value = '  -0.21%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.57%  '

# Row 19
$ws.Range("E19").Value = '  -3.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.32%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.76%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.737'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '

# Row 24
$ws.Range("E24").Value = '  -1.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '

# Row 26
$ws.Range("E26").Value = '  +4.74%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.35%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.44%  '

# Row 29
$ws.Range("E29").Value = '  +0.14%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.89%  '

# Row 31
$ws.Range("E31").Value = '  +1.61%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.51%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.90%  '

# Row 34
$ws.Range("E34").Value = '  -1.12%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '

# Row 36
$ws.Range("E36").Value = '  -1.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.26%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.140'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.68%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '482.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +14.12%  '

# Row 40
$ws.Range("E40").Value = '  +1.27%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.84%  '

# Row 42
$ws.Range("E42").Value = '  -2.14%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.12%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.81%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.82%  '

# Row 46
$ws.Range("D46").Value = '2.932.18'
$ws.Range("E46").Value = '  -2.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0362'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.62%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.41%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.99%  '
